$wb = $excel.ActiveWorkbook

# Sheet names affected: "展览" (Exhibition) and "全部类型" (All Types)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F6").Value  = 213
    $ws.Range("F10").Value = 5509
    $ws.Range("F11").Value = 11
    $ws.Range("F17").Value = 168
    $ws.Range("F18").Value = 220
}
